$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (this shifts existing D:K data to F:M,
# carrying over each columns number formatting/style along with it).
$ws.Range("D1:E1").EntireColumn.Insert()

# The two freshly inserted columns (D:E) do not yet carry the numeric/date styles
# that the rest of the quarter columns use, so copy formats across from the
# neighboring (just-shifted) F:M block, which is identical in shape to D:E+old cols.
$ws.Range("F7:M102").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new columns with the latest reported quarters financial data.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 60200
$ws.Range("E8").Value = 54900
$ws.Range("D9").Value = 22200
$ws.Range("E9").Value = 20100
$ws.Range("D10").Value = 38000
$ws.Range("E10").Value = 34800
$ws.Range("D12").Value = 12300
$ws.Range("E12").Value = 11700
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 73500
$ws.Range("E17").Value = 69800
$ws.Range("D18").Value = -13300
$ws.Range("E18").Value = -14900
$ws.Range("D20").Value = -500
$ws.Range("E20").Value = -100
$ws.Range("D21").Value = -13200
$ws.Range("E21").Value = -14500
$ws.Range("D22").Value = 100
$ws.Range("E22").Value = 100
$ws.Range("D23").Value = -13900
$ws.Range("E23").Value = -15100
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -13900
$ws.Range("E26").Value = -15000
$ws.Range("D27").Value = -13900
$ws.Range("E27").Value = -15000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 500
$ws.Range("E32").Value = 100
$ws.Range("D33").Value = -13900
$ws.Range("E33").Value = -15000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -13900
$ws.Range("E35").Value = -15000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 94900
$ws.Range("E41").Value = 107300
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 79400
$ws.Range("E43").Value = 62500
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 35300
$ws.Range("E45").Value = 18400
$ws.Range("D46").Value = 209600
$ws.Range("E46").Value = 188100
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 7500
$ws.Range("E48").Value = 3300
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 16000
$ws.Range("E52").Value = 15500
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 233200
$ws.Range("E54").Value = 206900
$ws.Range("D57").Value = 9200
$ws.Range("E57").Value = 6600
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 119200
$ws.Range("E59").Value = 104100
$ws.Range("D60").Value = 128400
$ws.Range("E60").Value = 110700
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 31600
$ws.Range("E62").Value = 14400
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 160000
$ws.Range("E66").Value = 125000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -145600
$ws.Range("E72").Value = -131700
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 73200
$ws.Range("E76").Value = 81800
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -13900
$ws.Range("E81").Value = -15000
$ws.Range("D83").Value = 600
$ws.Range("E83").Value = 500
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -7400
$ws.Range("E89").Value = -300
$ws.Range("D91").Value = -4800
$ws.Range("E91").Value = -600
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -4800
$ws.Range("E94").Value = -600
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 400
$ws.Range("E100").Value = 58500
$ws.Range("D101").Value = -600
$ws.Range("E101").Value = -600
$ws.Range("D102").Value = -12300
$ws.Range("E102").Value = 56900
